$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 234
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
